$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'332.01"
$ws.Range('E2').Value = "'0.92%"
$ws.Range('E3').Value = "'4.25%"
$ws.Range('D4').Value = "'5.628"
$ws.Range('E4').Value = "'2.57%"
$ws.Range('D5').Value = "'0.08356"
$ws.Range('E5').Value = "'4.44%"
$ws.Range('E6').Value = "'2.88%"
$ws.Range('B7').Value = 'GateToken'
$ws.Range('C7').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D7').Value = "'4.449"
$ws.Range('E7').Value = "'1.74%"
$ws.Range('B8').Value = 'MXToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D8').Value = "'0.9729"
$ws.Range('E8').Value = "'2.50%"
$ws.Range('B9').Value = 'BTSEToken'
$ws.Range('C9').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('D9').Value = "'2.583"
$ws.Range('E9').Value = "'0.34%"
$ws.Range('B10').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C10').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D10').Value = "'0.1165"
$ws.Range('E10').Value = "'3.45%"
$ws.Range('B11').Value = 'WazirX'
$ws.Range('C11').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D11').Value = "'0.1917"
$ws.Range('E11').Value = "'2.07%"
$ws.Range('B12').Value = 'MCDex'
$ws.Range('C12').Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range('D12').Value = "'10.36"
$ws.Range('E12').Value = "'-2.62%"
$ws.Range('B13').Value = 'MandalaExchangeToken'
$ws.Range('C13').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D13').Value = "'0.09951"
$ws.Range('E13').Value = "'-0.77%"
$ws.Range('B14').Value = 'BitrueCoin'
$ws.Range('C14').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D14').Value = "'0.04707"
$ws.Range('E14').Value = "'-1.96%"
$ws.Range('B15').Value = 'BitMartToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D15').Value = "'0.1059"
$ws.Range('E15').Value = "'0.15%"
$ws.Range('B16').Value = 'BitForexToken'
$ws.Range('C16').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D16').Value = "'0.001296"
$ws.Range('E16').Value = "'2.28%"
$ws.Range('B17').Value = 'TigerCash'
$ws.Range('C17').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D17').Value = "'0.006045"
$ws.Range('E17').Value = "'1.11%"
$ws.Range('B18').Value = 'LEO'
$ws.Range('C18').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D18').Value = "'3.372"
$ws.Range('E18').Value = "'0.23%"
$ws.Range('E19').Value = "'-3.06%"
$ws.Range('D20').Value = "'0.1391"
$ws.Range('E20').Value = "'-2.06%"
$ws.Range('D22').Value = "'0.04199"
$ws.Range('E22').Value = "'2.92%"
$ws.Range('D23').Value = "'0.001310"
$ws.Range('E23').Value = "'3.36%"
$ws.Range('D24').Value = "'0.004589"
$ws.Range('E24').Value = "'5.93%"
$ws.Range('D25').Value = "'0.0001303"
$ws.Range('E25').Value = "'8.53%"
$ws.Range('D26').Value = "'0.0003743"
$ws.Range('E26').Value = "'-0.04%"
$ws.Range('D38').Value = "'0.02762"
$ws.Range('E38').Value = "'6.91%"
$ws.Range('D39').Value = "'0.05810"
$ws.Range('E39').Value = "'2.94%"
$ws.Range('D40').Value = "'0.007679"
$ws.Range('E40').Value = "'1.96%"
$ws.Range('D41').Value = "'0.1436"
$ws.Range('E41').Value = "'2.99%"
$ws.Range('D42').Value = "'0.007278"
$ws.Range('E42').Value = "'-1.65%"
$ws.Range('D43').Value = "'0.002014"
$ws.Range('E43').Value = "'-0.07%"
$ws.Range('D44').Value = "'0.008041"
$ws.Range('E44').Value = "'-7.09%"
$ws.Range('D45').Value = "'0.3399"
$ws.Range('D46').Value = "'0.00007273"
$ws.Range('E46').Value = "'2.34%"
$ws.Range('D47').Value = "'0.00000000751"
$ws.Range('E47').Value = "'0.17%"
$ws.Range('D48').Value = "'0.0005805"
$ws.Range('E48').Value = "'-0.10%"
$ws.Range('D49').Value = "'0.003505"
$ws.Range('E49').Value = "'-3.91%"
$ws.Range('D50').Value = "'0.003501"
$ws.Range('E50').Value = "'-0.84%"
$ws.Range('D51').Value = "'0.00002104"
$ws.Range('E51').Value = "'0.17%"
